$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (also updates the <sheet name="..."> entry in workbook.xml)
$ws.Name = "Date1"

# Update row 2
$ws.Range("B2").Value = "1 bj3"
$ws.Range("C2").Value = "dfgbd"
$ws.Range("D2").Value = "fbv"

# Update row 3
$ws.Range("B3").Value = "2 bj3"
$ws.Range("C3").Value = "dfbv"
$ws.Range("D3").Value = "dfcv"

# Update row 4
$ws.Range("B4").Value = "3 bje"
$ws.Range("C4").Value = "sdzdgvc"
$ws.Range("D4").Value = "svc"
